$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("E2").Value = 23.19000000000019
$ws.Range("H2").Value = [double]"5.657187386625001e-16"
$ws.Range("K2").Value = 42.48955603030929
$ws.Range("L2").Value = "[33.003721140087954, 51.97539092053062]"
$ws.Range("O2").Value = 1.46544762419704
$ws.Range("P2").Value = "[1.2264475824825025, 1.7044476659115775]"
$ws.Range("S2").Value = 61.85006149620168
$ws.Range("T2").Value = "[55.992722676356735, 67.70740031604663]"
$ws.Range("W2").Value = 17.78132132132147
$ws.Range("X2").Value = 16.89921921921936
$ws.Range("Y2").Value = 18.66342342342357

# Row 3 updates
$ws.Range("B3").Value = 1
$ws.Range("E3").Value = 22.80000000000013
$ws.Range("G3").Value = [double]"1.641575764210756e-12"
$ws.Range("H3").Value = [double]"4.989947848860703e-12"
$ws.Range("K3").Value = 43.41263578142147
$ws.Range("L3").Value = "[29.122207634705696, 57.70306392813725]"
$ws.Range("M3").Value = [double]"1.025800910881003e-08"
$ws.Range("N3").Value = [double]"1.025800910881003e-08"
$ws.Range("O3").Value = 0.9622896416401163
$ws.Range("P3").Value = "[0.6226580034141929, 1.3019212798660398]"
$ws.Range("Q3").Value = [double]"7.886965924797096e-08"
$ws.Range("R3").Value = [double]"7.886965924797096e-08"
$ws.Range("S3").Value = 59.61818441733162
$ws.Range("T3").Value = "[51.96761067777891, 67.26875815688433]"
$ws.Range("W3").Value = 19.30810810810821
$ws.Range("X3").Value = 18.07567567567578
$ws.Range("Y3").Value = 20.54054054054065
